# Split the poster title run into three runs so the word "using" becomes
# "Using" (title-cased), matching the authored edit:
#   "Mapping Emotional Landscapes of Fiction using Machine Learning Techniques"
# -> "Mapping Emotional Landscapes of " + "Fiction Using " + "Machine Learning Techniques"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the title shape robustly by its current text rather than a hard-coded
# shape index.
$titleShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text.StartsWith("Mapping Emotional Landscapes")) {
                $titleShape = $shp
                break
            }
        }
    }
}

$titlePara = $titleShape.TextFrame.TextRange.Paragraphs(1, 1)

$firstPart  = "Mapping Emotional Landscapes of "
$middlePart = "Fiction Using "
$lastPart   = "Machine Learning Techniques"

# Replace each portion in place (right-to-left doesn't matter here since all
# three pieces together cover the whole original paragraph length), which
# splits the single run into three runs while preserving the existing
# character formatting (size 6000, gray fill, etc.) on each piece.
$titlePara.Characters(1, $firstPart.Length).Text = $firstPart
$titlePara.Characters($firstPart.Length + 1, $middlePart.Length).Text = $middlePart
$titlePara.Characters($firstPart.Length + $middlePart.Length + 1, $lastPart.Length).Text = $lastPart
